# Rework the workout plan: rename the sheet and swap in a new
# exercise/rep-scheme for every day of the program.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Week 1"

# DAY 1
$ws.Range("A7").Value  = "Incline Press Machine"
$ws.Range("D7").Value  = "3x8"
$ws.Range("A8").Value  = "Rear Delt Flies"
$ws.Range("D8").Value  = "4x12"
$ws.Range("A9").Value  = "DB Press"
$ws.Range("D9").Value  = "3x6"
$ws.Range("A10").Value = "Skullcrushers"
$ws.Range("D10").Value = "3x8"
$ws.Range("A11").Value = "Overhead Cable Extensions"
$ws.Range("D11").Value = "4x8"

# DAY 2
$ws.Range("A16").Value = "Seated Cable Rows"
$ws.Range("D16").Value = "4x12"
$ws.Range("A17").Value = "DB Rows"
$ws.Range("D17").Value = "4x12"
$ws.Range("A18").Value = "Straight Arm Lat Pulldowns"
$ws.Range("D18").Value = "4x10"
$ws.Range("A19").Value = "DB Curls"
$ws.Range("D19").Value = "3x6"
$ws.Range("A20").Value = "Waiter Curl"
$ws.Range("D20").Value = "4x10"

# DAY 3
$ws.Range("A25").Value = "Quad Extensions"
$ws.Range("D25").Value = "4x8"
$ws.Range("A26").Value = "Hack Squat"
$ws.Range("D26").Value = "3x6"
$ws.Range("A27").Value = "Split Squats"
$ws.Range("D27").Value = "4x6"
$ws.Range("A28").Value = "Romanian Deadlifts"
$ws.Range("D28").Value = "3x6"
$ws.Range("A29").Value = "Hip Thrusts"
$ws.Range("D29").Value = "3x6"

# DAY 4
$ws.Range("A34").Value = "Incline Press Machine"
$ws.Range("D34").Value = "4x6"
$ws.Range("A35").Value = "Lateral Raises"
$ws.Range("D35").Value = "4x10"
$ws.Range("A36").Value = "Cable Flies"
$ws.Range("D36").Value = "3x12"
$ws.Range("A37").Value = "Tricep Kickbacks"
$ws.Range("D37").Value = "3x6"
$ws.Range("A38").Value = "Rope Pushdowns"
$ws.Range("D38").Value = "3x12"

# DAY 5
$ws.Range("A43").Value = "Seated Cable Rows"
$ws.Range("D43").Value = "3x6"
$ws.Range("A44").Value = "Barbell Rows"
$ws.Range("D44").Value = "3x6"
$ws.Range("A45").Value = "Lat Pulldowns"
$ws.Range("D45").Value = "3x10"
$ws.Range("A46").Value = "Barbell Curls"
$ws.Range("D46").Value = "4x6"
$ws.Range("A47").Value = "Cable Curls"
$ws.Range("D47").Value = "4x8"

# DAY 6
$ws.Range("A52").Value = "Hip Thrusts"
$ws.Range("D52").Value = "3x8"
$ws.Range("A53").Value = "Goblet Squat"
$ws.Range("D53").Value = "3x10"
$ws.Range("A54").Value = "Romanian Deadlifts"
$ws.Range("D54").Value = "3x12"
$ws.Range("A55").Value = "Quad Extensions"
$ws.Range("D55").Value = "3x8"
$ws.Range("A56").Value = "Barbell Lunges"
$ws.Range("D56").Value = "3x10"
